$wb = $excel.ActiveWorkbook

$oldCreateName = "Create_Feature Category_x0009_"
$oldEditName = "Edit_Feature Category_x0009_"

# Add two brand-new sheets at the END of the workbook first. This advances the
# internal sheetId counter (old sheets are sheetId 1/2, so the new ones land on
# 3/4, matching the target) without disturbing the index positions of the two
# existing sheets (so we can keep re-fetching them safely by name).
$wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null

# The two new sheets are now the last two tabs (index 3 and 4).
$newCreate = $wb.Worksheets.Item(3)
$newEdit = $wb.Worksheets.Item(4)

# ---- "Create_Feature Category" sheet content (header row had its spaces
# stripped; data row is unchanged) ----
$newCreate.Range("A1").Value = "FeatureCategory"
$newCreate.Range("B1").Value = "CategoryDescription`t"
$newCreate.Range("C1").Value = "DefaultSelectionSequence`t"
$newCreate.Range("D1").Value = "FeatureID"
$newCreate.Range("E1").Value = "FeatureDescription"
$newCreate.Range("A2").Value = "PK_FC1"
$newCreate.Range("B2").Value = "Feature Category`t"
$newCreate.Range("C2").Value = 10
$newCreate.Range("D2").Value = "F1"
$newCreate.Range("E2").Value = "Color- Red"

# ---- "Edit_Feature Category" sheet content (header row had its spaces
# stripped; data row is unchanged) ----
$newEdit.Range("A1").Value = "FeatureCategory`t"
$newEdit.Range("B1").Value = "CategoryDescription`t"
$newEdit.Range("C1").Value = "DefaultSelectionSequence`t"
$newEdit.Range("A2").Value = "PK_FC1"
$newEdit.Range("B2").Value = "Feature Category - Update`t"
$newEdit.Range("C2").Value = 10

# Match the saved selections: sheet1 keeps the whole A1:E2 block selected,
# sheet2 ends up with B1 as the active/only selected cell.
$newCreate.Range("A1:E2").Select()
$newEdit.Range("B1").Select()

# Remove the original two sheets (looked up fresh by name right before the
# delete call, since sheet references track *index position*, not identity,
# across structural edits in this host).
$wb.Worksheets.Item($oldCreateName).Delete()
$wb.Worksheets.Item($oldEditName).Delete()

# The freshly-added sheets are now first and second; rename them to drop the
# stray trailing-tab (_x0009_) suffix and re-fetch fresh references.
$wb.Worksheets.Item(1).Name = "Create_Feature Category"
$wb.Worksheets.Item(2).Name = "Edit_Feature Category"

# "Edit_Feature Category" was the active tab before and stays the active tab.
$wb.Worksheets.Item(2).Activate()
